# Rename the worksheet from "Sheet1" to "replies"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "replies"

# --- Update / add cell values (Thought Record exercise extended with
#     find_automatic_thought, find_alternative_response and congratulations flows,
#     plus USER_NAME personalisation of several existing replies) ---
$ws.Cells.Item(1, 1).Value = 'intro'
$ws.Cells.Item(1, 2).Value = 'handle_sadness'
$ws.Cells.Item(1, 3).Value = 'goodbye'
$ws.Cells.Item(1, 4).Value = 'thought_record_intro'
$ws.Cells.Item(1, 5).Value = 'find_automatic_thought'
$ws.Cells.Item(1, 6).Value = 'thought_record_probing_questions'
$ws.Cells.Item(1, 7).Value = 'find_alternative_response'
$ws.Cells.Item(1, 8).Value = 'recommend_supervised_help'
$ws.Cells.Item(1, 9).Value = 'congratulations'
$ws.Cells.Item(2, 1).Value = 'Hi, I''m CloudBot 😊'
$ws.Cells.Item(2, 2).Value = 'I am here for you, USER_NAME. Together we can pass over it, ok?'
$ws.Cells.Item(2, 3).Value = 'Ok, USER_NAME, I am more than happy to see that you found something that can improve your mood ❤'
$ws.Cells.Item(2, 4).Value = 'Now, I would like to do a little exercise with you, USER_NAME.'
$ws.Cells.Item(2, 5).Value = 'Ok, USER_NAME, now that you completed the steps from 1 to 3, let''s move forward! ✨'
$ws.Cells.Item(2, 6).Value = 'Question 1: What is the effect of believing this thought?'
$ws.Cells.Item(2, 7).Value = 'Now, USER_NAME, use your responses to these questions to create an alternative response'
$ws.Cells.Item(2, 8).Value = 'It seems that this solution doesn''t help you enough, USER_NAME…'
$ws.Cells.Item(2, 9).Value = 'Good job, USER_NAME! 🎉'
$ws.Cells.Item(3, 1).Value = 'Let me introduce myself. I am a trained medical bot used mostly for managing mental health issues using Cognitive Behavioral Therapy (CBT) techniques.'
$ws.Cells.Item(3, 2).Value = 'I assume that this feeling can impact your daily activities 🥺'
$ws.Cells.Item(3, 3).Value = 'If you need something, do not forget to contact me or our amazing team!'
$ws.Cells.Item(3, 4).Value = 'In medical terms, it is called Thought Record.'
$ws.Cells.Item(3, 5).Value = 'Please pick one automatic thought from your list'
$ws.Cells.Item(3, 6).Value = 'Question 2: What would happen if you didn’t believe this thought?'
$ws.Cells.Item(3, 7).Value = 'This answer could be used for defending the automatic thought that may cause your negative emotion'
$ws.Cells.Item(3, 8).Value = 'I am sorry to see that…'
$ws.Cells.Item(3, 9).Value = 'I am more than happy to see that you succeed in reconstructing the way in which you think!'
$ws.Cells.Item(4, 1).Value = 'If you want to learn more about these techniques, I encourage you to visit the link below '
$ws.Cells.Item(4, 2).Value = 'But listen, everything that you need to cope with this sentiment is within you.'
$ws.Cells.Item(4, 3).Value = 'We are available at any time of the day.'
$ws.Cells.Item(4, 4).Value = 'The skill we will be working on here is called cognitive restructuring.'
$ws.Cells.Item(4, 5).Value = 'Considering it, you should respond to the following questions'
$ws.Cells.Item(4, 6).Value = 'Question 3: What is the evidence supporting this thought?'
$ws.Cells.Item(4, 7).Value = 'When you are ready, please let me know what idea you developed 🤗'
$ws.Cells.Item(4, 8).Value = 'But don’t worry, our amazing team is here to help you find a better solution for your problem ❤️'
$ws.Cells.Item(4, 9).Value = 'Don''t forget to apply this exercise every time to feel overwhelmed by a thought.'
$ws.Cells.Item(5, 1).Value = 'https://cogbtherapy.com/cognitive-behavior-therapy-techniques'
$ws.Cells.Item(5, 3).Value = 'Until next time, take care of yourself 🤗'
$ws.Cells.Item(5, 4).Value = 'Cognitive restructuring refers to the act of identifying ineffective patterns in thinking, and changing them to be more effective.'
$ws.Cells.Item(5, 5).Value = 'So, USER_NAME, what automatic thought bothers you the most? 💫'
$ws.Cells.Item(5, 6).Value = 'Question 4: What is the evidence against this thought?'
$ws.Cells.Item(5, 8).Value = 'They are all more than happy to help you out, so don''t hesitate to contact them!'
$ws.Cells.Item(6, 1).Value = 'I am here to help you whenever you may feel down or you may need somebody to talk to.'
$ws.Cells.Item(6, 3).Value = 'Have a wonderful day, USER_NAME! 🌸'
$ws.Cells.Item(6, 4).Value = 'More effective can mean triggering less negative emotion, seeing things more clearly, or enabling more skillful behavior. '
$ws.Cells.Item(6, 6).Value = 'Question 5: What’s the worst that could happen, and would you survive it?'
$ws.Cells.Item(6, 8).Value = 'Nevertheless, I am always here in case you need to talk to somebody'
$ws.Cells.Item(7, 1).Value = 'Now, what about you? What''s your name?'
$ws.Cells.Item(7, 4).Value = 'Cognitive restructuring builds on your ability to accurately recognize automatic thoughts and feelings. '
$ws.Cells.Item(7, 6).Value = 'Question 6: What’s the best that could happen?'
$ws.Cells.Item(7, 8).Value = 'You will be better soon, I promise '
$ws.Cells.Item(8, 4).Value = 'Sounds interesting, right? 😊'
$ws.Cells.Item(8, 6).Value = 'Question 7: What is the most likely?'
$ws.Cells.Item(8, 8).Value = 'Until next time, take care of yourself 🤗'
$ws.Cells.Item(9, 4).Value = 'Then let''s see the steps that we have to check in order to achieve this new skill, USER_NAME! ✨'
$ws.Cells.Item(9, 6).Value = 'Question 8: If your friend was in this situation, what would you tell him/her?'
$ws.Cells.Item(9, 8).Value = 'Have a wonderful day, USER_NAME! 🌸'
$ws.Cells.Item(10, 4).Value = 'As a little tip, you might feel the need to grab a pen and a piece of paper near you'
$ws.Cells.Item(10, 6).Value = 'Question 9: What can you do about this?'
$ws.Cells.Item(11, 4).Value = 'Step 1: Identify the situation that triggered the negative emotion that you feel right now 🌪️'
$ws.Cells.Item(12, 4).Value = 'Step 2: Write down the emotions that triggered because of the situation. Rate their intensity from 0 to 10 💫'
$ws.Cells.Item(13, 4).Value = 'Step 3: List all of your automatic thoughts that came in your mind and rate how much you believe them on a scale from 0-10 ☘️'
$ws.Cells.Item(14, 4).Value = 'Step 4: Choose the automatic thought that is most responsible for your distress, and use it to answer the probing questions 🌤️'
$ws.Cells.Item(15, 4).Value = 'Step 5: Using your answers to the probing questions, develop a short alternative response to the automatic thought you choose 🌞'
$ws.Cells.Item(16, 4).Value = 'Now I will let you do the steps from 1 to 3, USER_NAME.'
$ws.Cells.Item(17, 4).Value = 'In general, it takes between 1 and 3 minutes to complete them, but I encourage you to take your time and reflect ✨'
$ws.Cells.Item(18, 4).Value = 'I know that this activity can be emotionally consuming, so just type READY when you finish 😊'

# --- Clear the cells that no longer hold data in column E (data moved to column F) ---
$ws.Range("E6:E10").ClearContents()

# --- Re-apply the (empty, wrap-text styled) placeholder cells E13:E14 next to the
#     Step 3 / Step 4 instructions, matching D13:D14 formatting ---
$ws.Cells.Item(13, 5).WrapText = $true
$ws.Cells.Item(14, 5).WrapText = $true

# --- Column widths: column D (Thought Record) is duplicated into the new column E
#     (find_automatic_thought), old column E (probing questions) shifts to F, and three
#     new columns G/H/I are added for the new conversation flows ---
$ws.Columns.Item(5).ColumnWidth = 104
$ws.Columns.Item(6).ColumnWidth = 58.833333333333336
$ws.Columns.Item(7).ColumnWidth = 72.33333333333333
$ws.Columns.Item(8).ColumnWidth = 76.83333333333333
$ws.Columns.Item(9).ColumnWidth = 72.5

# --- Restore the view: scrolled to column D, active cell E5 ---
$ws.Range("E5").Select()
